$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy the date format already used in A5:A7 onto the new date cells, then
# set the values/contents for the two new benchmark rows.
$ws.Range("A5").Copy()
$ws.Range("A8:A9").PasteSpecial(-4122)

# Row 8: A bit better.
$ws.Range("A8").Value = 44498
$ws.Range("B8").Value = "7.0.1 (develop)"
$ws.Range("C8").Value = 589.29999999999995
$ws.Range("D8").Value = 20.829999999999998
$ws.Range("E8").Value = 262
$ws.Range("F8").Value = 232.59
$ws.Range("H8").Value = "Intel core i9-9900K@3.6GHz"
$ws.Range("L8").Value = "A bit better."

# Row 9: Removed excel checking
$ws.Range("A9").Value = 44498
$ws.Range("B9").Value = "7.0.1 (develop)"
$ws.Range("C9").Value = 469.01999999999998
$ws.Range("D9").Value = 19.469999999999999
$ws.Range("E9").Value = 453
$ws.Range("F9").Value = 107.65000000000001
$ws.Range("H9").Value = "Intel core i9-9900K@3.6GHz"
$ws.Range("L9").Value = "Removed excel checking"

$ws.Range("H9").Select() | Out-Null
